# CORRIDAS TK MULTIMARCAS.xlsx update
# - "MES 01": a few separator-row date values advance by one day, a
#   previously-empty data row gets a new destination ("SONHOS DOURADOS"),
#   the trailing two rows (old 66/67) are removed, the TOTAL row's SUM
#   formula drops the now-gone A67 term, and the date-separator rows swap
#   which border/format variant they use.
# - Window / view geometry is refreshed to match where the author left the
#   workbook scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")

# --- Value edits in the tail of the sheet -------------------------------
$ws.Cells.Item(60, 1).Value2 = 45218
$ws.Cells.Item(61, 1).Value2 = 15
$ws.Cells.Item(61, 2).Value2 = "SONHOS DOURADOS"
$ws.Cells.Item(62, 1).Value2 = 45219
$ws.Cells.Item(64, 1).Value2 = 45220

# --- Swap the date-separator row formatting between the two variants ----
# Stash the two existing "template" formats (rows 2 and 16) in scratch
# cells far off to the side before overwriting either of them.
$ws.Range("A2:B2").Copy()
$ws.Range("Z1:AA1").PasteSpecial(-4122)

$ws.Range("A16:B16").Copy()
$ws.Range("Z2:AA2").PasteSpecial(-4122)

$rowsToVariantB = @(2, 4, 6, 8, 14)
foreach ($r in $rowsToVariantB) {
  $ws.Range("Z2:AA2").Copy()
  $ws.Range("A$r`:B$r").PasteSpecial(-4122)
}

$rowsToVariantA = @(16, 21, 26, 29, 32, 40, 44, 48, 53, 55, 58, 60, 62, 64)
foreach ($r in $rowsToVariantA) {
  $ws.Range("Z1:AA1").Copy()
  $ws.Range("A$r`:B$r").PasteSpecial(-4122)
}

$ws.Range("Z1:AA2").Clear()

# --- Drop the now-unused trailing rows and fix the TOTAL formula --------
$ws.Rows.Item(67).Delete()
$ws.Rows.Item(66).Delete()
$ws.Cells.Item(66, 2).Formula = "=SUM(A3,A5,A7,A9:A13,A15,A17:A20,A22:A25,A27:A28,A30:A31,A33:A39,A41:A43,A45:A47,A49:A52,A54,A56:A57,A59,A61,A63,A65)"

# --- Refresh the view so it lands where the author left it --------------
$ws.Application.ActiveWindow.ScrollRow = 45
$ws.Range("B67").Select()

$excel.ActiveWindow.WindowState = -4143
